$wb = $excel.ActiveWorkbook

# Sheet ALC, row 6 (diff hunk @@ -932,22 +932,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 49.8
$ws.Range("I6").Value = 49.75
$ws.Range("K6").Value = 149.25
$ws.Range("M6").Value = -37.25

# Sheet ALC, row 40 (diff hunk @@ -2604,25 +2604,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3782.8572
$ws.Range("J40").Value = 2590
$ws.Range("L40").Value = 2590
$ws.Range("N40").Value = -2940

# Sheet ALC, row 112 (diff hunk @@ -6186,25 +6186,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3675.4
$ws.Range("J112").Value = 2933.2104
$ws.Range("L112").Value = 8799.6312
$ws.Range("N112").Value = -11015.6312

# Sheet ALC, row 113 (diff hunk @@ -6238,25 +6238,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3461.4119
$ws.Range("I113").Value = 2693.4285
$ws.Range("J113").Value = 3999
$ws.Range("K113").Value = 2693.4285
$ws.Range("L113").Value = 3999
$ws.Range("M113").Value = 560.5715
$ws.Range("N113").Value = -10507

# Sheet ALC, row 138 (diff hunk @@ -7469,25 +7469,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6520.25
$ws.Range("J138").Value = 6944.905
$ws.Range("L138").Value = 20834.715
$ws.Range("N138").Value = -31114.715

# Sheet ARM, row 61 (diff hunk @@ -10644,25 +10644,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3901.8647
$ws.Range("I61").Value = 3889.6333
$ws.Range("J61").Value = 3954.2856
$ws.Range("K61").Value = 3889.6333
$ws.Range("L61").Value = 3954.2856
$ws.Range("M61").Value = -3677.6333
$ws.Range("N61").Value = -4378.2856

# Sheet ARM, row 132 (diff hunk @@ -14048,22 +14048,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1523978.1
$ws.Range("I132").Value = 2268543.5
$ws.Range("K132").Value = 6805630.5
$ws.Range("M132").Value = -6803100.5

# Sheet ARM, row 136 (diff hunk @@ -14247,25 +14247,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3901.8647
$ws.Range("I136").Value = 3889.6333
$ws.Range("J136").Value = 3954.2856
$ws.Range("K136").Value = 11668.8999
$ws.Range("L136").Value = 11862.8568
$ws.Range("M136").Value = -9118.8999
$ws.Range("N136").Value = -16962.8568

# Sheet BSM, row 38 (diff hunk @@ -16378,20 +16378,23 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 26530
$ws.Range("I38").Value = 12000
$ws.Range("K38").Value = 12000
$ws.Range("M38").Value = -11584

# Sheet BSM, row 107 (diff hunk @@ -19711,25 +19714,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2942776.8
$ws.Range("I107").Value = 4168390
$ws.Range("J107").Value = 1305.2
$ws.Range("K107").Value = 4168390
$ws.Range("L107").Value = 1305.2
$ws.Range("M107").Value = -4166470
$ws.Range("N107").Value = -5145.2

# Sheet BSM, row 109 (diff hunk @@ -19809,22 +19812,19 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Sheet CRP, row 10 (diff hunk @@ -21867,22 +21867,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1999.5
$ws.Range("I10").Value = 1999
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 1999
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = -1860
$ws.Range("N10").Value = -2278

# Sheet CRP, row 14 (diff hunk @@ -22060,22 +22063,19 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# Sheet CRP, row 15 (diff hunk @@ -22109,25 +22109,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1500
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 1500
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 1500
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -1840

# Sheet CRP, row 31 (diff hunk @@ -22899,22 +22896,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4717.36
$ws.Range("I31").Value = 1300.2858
$ws.Range("K31").Value = 1300.2858
$ws.Range("M31").Value = -1005.2858

# Sheet CRP, row 34 (diff hunk @@ -23049,22 +23046,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4717.36
$ws.Range("I34").Value = 1300.2858
$ws.Range("K34").Value = 1300.2858
$ws.Range("M34").Value = -1098.2858

# Sheet CRP, row 107 (diff hunk @@ -26596,22 +26593,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 756.4
$ws.Range("I107").Value = 555.1429000000001
$ws.Range("K107").Value = 555.1429000000001
$ws.Range("M107").Value = 1364.8571

# Sheet CRP, row 121 (diff hunk @@ -27258,19 +27255,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H121").Value = 49000
$ws.Range("J121").Value = 49000
$ws.Range("L121").Value = 49000
$ws.Range("N121").Value = -51620

# Sheet CRP, row 132 (diff hunk @@ -27782,22 +27782,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4057.3845
$ws.Range("I132").Value = 3241.4783
$ws.Range("K132").Value = 9724.4349
$ws.Range("M132").Value = -7194.4349

# Sheet CRP, row 133 (diff hunk @@ -27834,19 +27834,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060

# Sheet CUL, row 8 (diff hunk @@ -28678,22 +28681,22 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1372.7142
$ws.Range("I8").Value = 1372.7142
$ws.Range("K8").Value = 4118.142599999999
$ws.Range("M8").Value = -3979.142599999999

# Sheet CUL, row 122 (diff hunk @@ -34456,25 +34459,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 124977.39
$ws.Range("I122").Value = 313.2
$ws.Range("J122").Value = 165191.64
$ws.Range("K122").Value = 2818.8
$ws.Range("L122").Value = 1486724.76
$ws.Range("M122").Value = -368.7999999999997
$ws.Range("N122").Value = -1491624.76

# Sheet GSM, row 5 (diff hunk @@ -35722,22 +35725,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 300
$ws.Range("I5").Value = 300
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 300
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -188
$ws.Range("N5").ClearContents()

# Sheet GSM, row 9 (diff hunk @@ -35918,25 +35921,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1136
$ws.Range("J9").Value = 272.5
$ws.Range("L9").Value = 272.5
$ws.Range("N9").Value = -612.5

# Sheet GSM, row 22 (diff hunk @@ -36558,22 +36561,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 32000
$ws.Range("I22").Value = 50000
$ws.Range("J22").Value = 14000
$ws.Range("K22").Value = 50000
$ws.Range("L22").Value = 14000
$ws.Range("M22").Value = -49471
$ws.Range("N22").Value = -15058

# Sheet GSM, row 107 (diff hunk @@ -40663,25 +40669,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 854.7143
$ws.Range("J107").Value = 940
$ws.Range("L107").Value = 940
$ws.Range("N107").Value = -4780

# Sheet GSM, row 132 (diff hunk @@ -41849,25 +41855,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 45458424
$ws.Range("I132").Value = 66670384
$ws.Range("J132").Value = 4220.143
$ws.Range("K132").Value = 200011152
$ws.Range("L132").Value = 12660.429
$ws.Range("M132").Value = -200008622
$ws.Range("N132").Value = -17720.429

# Sheet LTW, row 4 (diff hunk @@ -42525,22 +42531,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 25012334
$ws.Range("I4").Value = 25011250
$ws.Range("J4").Value = 25014500
$ws.Range("K4").Value = 25011250
$ws.Range("L4").Value = 25014500
$ws.Range("M4").Value = -25011137
$ws.Range("N4").Value = -25014726

# Sheet LTW, row 5 (diff hunk @@ -42574,19 +42583,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 29428.285
$ws.Range("I5").Value = 33333
$ws.Range("J5").Value = 6000
$ws.Range("K5").Value = 33333
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = -33220
$ws.Range("N5").Value = -6226

# Sheet LTW, row 7 (diff hunk @@ -42666,25 +42681,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5071.4253
$ws.Range("I7").Value = 4781.0835
$ws.Range("J7").Value = 5374.391
$ws.Range("K7").Value = 4781.0835
$ws.Range("L7").Value = 5374.391
$ws.Range("M7").Value = -4669.0835
$ws.Range("N7").Value = -5598.391

# Sheet LTW, row 10 (diff hunk @@ -42816,25 +42831,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 14290273
$ws.Range("I10").Value = 16669485
$ws.Range("J10").Value = 15000
$ws.Range("K10").Value = 16669485
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = -16669345
$ws.Range("N10").Value = -15280

# Sheet LTW, row 21 (diff hunk @@ -43361,19 +43376,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 13400
$ws.Range("J21").Value = 13400
$ws.Range("L21").Value = 13400
$ws.Range("N21").Value = -13748

# Sheet LTW, row 22 (diff hunk @@ -43407,25 +43425,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5010.724
$ws.Range("I22").Value = 3298.4348
$ws.Range("J22").Value = 11574.5
$ws.Range("K22").Value = 3298.4348
$ws.Range("L22").Value = 11574.5
$ws.Range("M22").Value = -3003.4348
$ws.Range("N22").Value = -12164.5

# Sheet LTW, row 27 (diff hunk @@ -43661,25 +43679,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 5010.724
$ws.Range("I27").Value = 3298.4348
$ws.Range("J27").Value = 11574.5
$ws.Range("K27").Value = 3298.4348
$ws.Range("L27").Value = 11574.5
$ws.Range("M27").Value = -3191.4348
$ws.Range("N27").Value = -11788.5

# Sheet LTW, row 28 (diff hunk @@ -43713,22 +43731,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H28").Value = 25012334
$ws.Range("I28").Value = 25011250
$ws.Range("J28").Value = 25014500
$ws.Range("K28").Value = 25011250
$ws.Range("L28").Value = 25014500
$ws.Range("M28").Value = -25011018
$ws.Range("N28").Value = -25014964

# Sheet LTW, row 37 (diff hunk @@ -44154,22 +44175,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H37").Value = 25012334
$ws.Range("I37").Value = 25011250
$ws.Range("J37").Value = 25014500
$ws.Range("K37").Value = 25011250
$ws.Range("L37").Value = 25014500
$ws.Range("M37").Value = -25011143
$ws.Range("N37").Value = -25014714

# Sheet LTW, row 126 (diff hunk @@ -48476,25 +48500,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5071.4253
$ws.Range("I126").Value = 4781.0835
$ws.Range("J126").Value = 5374.391
$ws.Range("K126").Value = 14343.2505
$ws.Range("L126").Value = 16123.173
$ws.Range("M126").Value = -11873.2505
$ws.Range("N126").Value = -21063.173

# Sheet WVR, row 8 (diff hunk @@ -49636,22 +49660,22 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 50000000
$ws.Range("I8").Value = 50000000
$ws.Range("K8").Value = 50000000
$ws.Range("M8").Value = -49999860

# Sheet WVR, row 29 (diff hunk @@ -50659,25 +50683,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 23399.8
$ws.Range("I29").Value = 33333
$ws.Range("J29").Value = 8500
$ws.Range("K29").Value = 33333
$ws.Range("L29").Value = 8500
$ws.Range("M29").Value = -33043
$ws.Range("N29").Value = -9080
